$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 79; existing rows 79-89 shift down to 80-90.
$ws.Rows(79).Insert()

# Populate the new row 79 with the new weekly price-report entry.
$ws.Range("A79").Value2 = 10
$ws.Range("B79").Value2 = "Vega Modelo de Temuco"
$ws.Range("C79").Value2 = "La Araucanía"
$ws.Range("D79").Value2 = 45124
$ws.Range("E79").Value2 = 9
$ws.Range("F79").Value2 = 100112042
$ws.Range("G79").Value2 = "Locoto"
$ws.Range("H79").Value2 = "Sin especificar"
$ws.Range("I79").Value2 = "Primera"
$ws.Range("J79").Value2 = 150
$ws.Range("K79").Value2 = 3800
$ws.Range("L79").Value2 = 3800
$ws.Range("M79").Value2 = 3800
$ws.Range("N79").Value2 = "`$/kilo"
$ws.Range("O79").Value2 = "Región de Arica y Parinacota"
$ws.Range("P79").Value2 = 3800
$ws.Range("Q79").Value2 = 1
$ws.Range("R79").Value2 = "Hortaliza"
